$d = $word.ActiveDocument

# Locate the anchor paragraph: "Java pattern matching statements case classes. ..."
$anchorText = "Java pattern matching statements case classes. Resource Monad wrapping Case classes."
$findRange = $d.Content
[void]$findRange.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Determine the 1-based index (within $d.Paragraphs) of the paragraph that
# contains the found text, by matching the paragraph whose range covers the
# start of the found range.
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  if ($p.Range.Start -le $findRange.Start -and $findRange.Start -lt $p.Range.End) {
    $anchorIndex = $i
    break
  }
}

$r = $d.Paragraphs.Item($anchorIndex).Range.Duplicate
$r.Collapse(0)   # wdCollapseEnd

# New list items to add; each inherits the anchor paragraph's list formatting
# (numId=3, ilvl=0, ind left=600 hanging=360) because InsertParagraphAfter()
# splits off of the existing list paragraph.
$texts = @(
  "Inputs / Sync Adapters:",
  "(Class, Instance, Member, Value) Events.",
  "Data Modelling Resources / Patterns:",
  "(Metaclass, Class, Instance, Context, Role, Occurrence);",
  "Switch actions: Populate Models (RDFS, OWL, Sets, Functional MVC / DCI DOM).",
  "Patterns Layers / Resource Monads hierarchies."
)

$insertAfterIndex = $anchorIndex
foreach ($t in $texts) {
  $insertAfterIndex = $insertAfterIndex + 1
  [void]$r.InsertParagraphAfter()
  $newPar = $d.Paragraphs.Item($insertAfterIndex)
  $newPar.Range.Text = $t
  $r = $newPar.Range.Duplicate
  $r.Collapse(0)
}
